# Insert a new data row before the current row 126 (Macroferia Regional de
# Talca - Betarraga). This shifts the existing rows 126..259 down to
# 127..260 and grows the sheet dimension to R260.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with its data values. The
# non-varying columns mirror every other data row in this sheet.
$ws.Cells.Item(126, 1).Value = 5
$ws.Cells.Item(126, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(126, 3).Value = "Maule"
$ws.Cells.Item(126, 4).Value = 44601
$ws.Cells.Item(126, 5).Value = 7
$ws.Cells.Item(126, 6).Value = 100114014
$ws.Cells.Item(126, 7).Value = "Betarraga"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Segunda"
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 700
$ws.Cells.Item(126, 12).Value = 700
$ws.Cells.Item(126, 13).Value = 700
$ws.Cells.Item(126, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(126, 15).Value = "Región del Maule"
$ws.Cells.Item(126, 16).Value = 140
$ws.Cells.Item(126, 17).Value = 5
$ws.Cells.Item(126, 18).Value = "Hortaliza"
